# Update 想去人数 (F) and one 最低票价 (G) value across all 4 sheets
# per commit: Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1752
$ws.Range("F4").Value = 169
$ws.Range("F5").Value = 478
$ws.Range("F6").Value = 839
$ws.Range("F7").Value = 261
$ws.Range("F8").Value = 1234
$ws.Range("F9").Value = 358
$ws.Range("F11").Value = 890
$ws.Range("F13").Value = 193
$ws.Range("F14").Value = 524
$ws.Range("F18").Value = 2969
$ws.Range("F19").Value = 2633
$ws.Range("F23").Value = 316
$ws.Range("F24").Value = 235
$ws.Range("F26").Value = 5334
$ws.Range("F29").Value = 27
$ws.Range("F30").Value = 60
$ws.Range("F31").Value = 335
$ws.Range("F35").Value = 296

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1143
$ws.Range("F7").Value = 233
$ws.Range("F8").Value = 20
$ws.Range("F24").Value = 321
$ws.Range("F25").Value = 282
$ws.Range("F26").Value = 3971
$ws.Range("F31").Value = 54
$ws.Range("G31").Value = 380
$ws.Range("F34").Value = 35
$ws.Range("F35").Value = 12
$ws.Range("F36").Value = 9

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2492
$ws.Range("F6").Value = 1060
$ws.Range("F9").Value = 1347
$ws.Range("F10").Value = 369
$ws.Range("F11").Value = 103

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2492
$ws.Range("F6").Value = 1752
$ws.Range("F7").Value = 1060
$ws.Range("F8").Value = 1347
$ws.Range("F9").Value = 369
$ws.Range("F10").Value = 103
$ws.Range("F11").Value = 169
$ws.Range("F12").Value = 478
$ws.Range("F13").Value = 839
$ws.Range("F14").Value = 261
$ws.Range("F15").Value = 1234
$ws.Range("F16").Value = 358
$ws.Range("F17").Value = 890
$ws.Range("F19").Value = 1143
$ws.Range("F20").Value = 1143
$ws.Range("F21").Value = 193
$ws.Range("F22").Value = 524
$ws.Range("F25").Value = 2969
$ws.Range("F26").Value = 2633
$ws.Range("F28").Value = 316
$ws.Range("F30").Value = 235
$ws.Range("F32").Value = 5334
$ws.Range("F37").Value = 27
$ws.Range("F38").Value = 60
$ws.Range("F39").Value = 335
$ws.Range("F44").Value = 321
$ws.Range("F45").Value = 321
$ws.Range("F48").Value = 54
$ws.Range("F51").Value = 296
